$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# "Elimina EC anteriores y se agregan nuevos": the list of mora periods
# (rows 16-22, column E) is regenerated in reverse order, and the two
# "Valor Mora" amounts that differed (column F) swap rows along with it.
$ws.Range("E16").Value = "2108"
$ws.Range("E17").Value = "2107"
$ws.Range("E18").Value = "2106"
$ws.Range("E19").Value = "2105"
$ws.Range("E20").Value = "2104"
$ws.Range("E21").Value = "2103"
$ws.Range("E22").Value = "2102"

$ws.Range("F16").Value = 31495
$ws.Range("F22").Value = 36341
